# Scheduled runner refresh: updates currentAveragePrice*/LevePrice*/LeveProfit*
# columns (H:N) with freshly-pulled market data for the affected leve rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row92 - Whinier than the Sword / Enchanted Koppranickel Ink
$ws.Range("H92").Value = 418.53333
$ws.Range("I92").Value = 362.64
$ws.Range("J92").Value = 698
$ws.Range("K92").Value = 362.64
$ws.Range("L92").Value = 698
$ws.Range("M92").Value = 885.36
$ws.Range("N92").Value = -3194

# ALC!row98 - The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 5395.6665
$ws.Range("I98").Value = 4687.3076
$ws.Range("J98").Value = 10000
$ws.Range("K98").Value = 4687.3076
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = -3189.3076
$ws.Range("N98").Value = -12996

# ALC!row122 - Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 5395.6665
$ws.Range("I122").Value = 4687.3076
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 14061.9228
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -11611.9228
$ws.Range("N122").Value = -34900

# ALC!row125 - Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 1419.579
$ws.Range("I125").Value = 985
$ws.Range("K125").Value = 8865
$ws.Range("M125").Value = -6405

# ALC!row129 - Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 6411512.5
$ws.Range("I129").Value = 83335570
$ws.Range("J129").Value = 1174.6945
$ws.Range("K129").Value = 250006710
$ws.Range("L129").Value = 3524.0835
$ws.Range("M129").Value = -250001710
$ws.Range("N129").Value = -13524.0835

$ws = $wb.Worksheets.Item("ARM")
# ARM!row3 - Skillet Labor / Bronze Skillet
$ws.Range("H3").Value = 3042.5
$ws.Range("I3").Value = 3042.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3042.5
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = -2927.5
$ws.Range("M3").ClearContents()

# ARM!row32 - Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 9046.298000000001
$ws.Range("I32").Value = 8187.793
$ws.Range("J32").Value = 19716.285
$ws.Range("K32").Value = 8187.793
$ws.Range("L32").Value = 19716.285
$ws.Range("M32").Value = -7900.793
$ws.Range("N32").Value = -20290.285

# ARM!row74 - As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 1639.1389
$ws.Range("I74").Value = 1147.5
$ws.Range("K74").Value = 1147.5
$ws.Range("M74").Value = -273.5

# ARM!row77 - Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 1639.1389
$ws.Range("I77").Value = 1147.5
$ws.Range("K77").Value = 5737.5
$ws.Range("M77").Value = -1369.5

$ws = $wb.Worksheets.Item("BSM")
# BSM!row23 - Get a Little Bit Closer / Brass Knuckles
$ws.Range("H23").Value = 58855.125
$ws.Range("J23").Value = 84008.2
$ws.Range("L23").Value = 84008.2
$ws.Range("N23").Value = -84574.2

# BSM!row31 - When Rhalgr Met Nophica / Spiked Knuckles
$ws.Range("H31").Value = 51007
$ws.Range("I31").Value = 38000
$ws.Range("J31").Value = 52865.145
$ws.Range("K31").Value = 38000
$ws.Range("L31").Value = 52865.145
$ws.Range("M31").Value = -37748
$ws.Range("N31").Value = -53369.145

# BSM!row94 - High Steal / High Steel Nugget
$ws.Range("H94").Value = 938.6667
$ws.Range("I94").Value = 456.2857
$ws.Range("J94").Value = 1360.75
$ws.Range("K94").Value = 456.2857
$ws.Range("L94").Value = 1360.75
$ws.Range("M94").Value = -5.28570000000002
$ws.Range("N94").Value = -2262.75

$ws = $wb.Worksheets.Item("CRP")
# CRP!row12 - A Sword in Hand / Ash Macuahuitl
$ws.Range("H12").Value = 6693082.5
$ws.Range("I12").Value = 10000621
$ws.Range("J12").Value = 78004.5
$ws.Range("K12").Value = 10000621
$ws.Range("L12").Value = 78004.5
$ws.Range("M12").Value = -10000451
$ws.Range("N12").Value = -78344.5

$ws = $wb.Worksheets.Item("CUL")
# CUL!row10 - A Real Fungi / Chanterelle Saute
$ws.Range("H10").Value = 450
$ws.Range("I10").Value = 42.5
$ws.Range("K10").Value = 127.5
$ws.Range("M10").Value = 11.5

# CUL!row113 - Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 2326936
$ws.Range("I113").Value = 5557606
$ws.Range("J113").Value = 853.6
$ws.Range("K113").Value = 16672818
$ws.Range("L113").Value = 2560.8
$ws.Range("M113").Value = -16670648
$ws.Range("N113").Value = -6900.8

# CUL!row131 - The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 1096.8276
$ws.Range("I131").Value = 1552.9
$ws.Range("J131").Value = 1037.5974
$ws.Range("K131").Value = 4658.700000000001
$ws.Range("L131").Value = 3112.7922
$ws.Range("M131").Value = 381.2999999999993
$ws.Range("N131").Value = -13192.7922

$ws = $wb.Worksheets.Item("GSM")
# GSM!row12 - Horn of Plenty / Bone Armillae
$ws.Range("H12").Value = 3372040.5
$ws.Range("I12").Value = 2882411.8
$ws.Range("K12").Value = 2882411.8
$ws.Range("M12").Value = -2882271.8

# GSM!row18 - Gorgeous Gorget / Brass Gorget
$ws.Range("H18").Value = 217148
$ws.Range("I18").Value = 1000000
$ws.Range("K18").Value = 1000000
$ws.Range("M18").Value = -999707

$ws = $wb.Worksheets.Item("LTW")
# LTW!row7 - Tan Before the Ban / Leather
$ws.Range("H7").Value = 2109.0908
$ws.Range("I7").Value = 2316
$ws.Range("J7").Value = 1936.6666
$ws.Range("K7").Value = 2316
$ws.Range("L7").Value = 1936.6666
$ws.Range("M7").Value = -2204
$ws.Range("N7").Value = -2160.6666

# LTW!row9 - From the Sands to the Stage / Leather Himantes
$ws.Range("H9").Value = 1042.909
$ws.Range("I9").Value = 271.42856
$ws.Range("J9").Value = 2393
$ws.Range("K9").Value = 271.42856
$ws.Range("L9").Value = 2393
$ws.Range("M9").Value = -47.42856
$ws.Range("N9").Value = -2841

# LTW!row43 - Subordinate Clause / Goatskin Choker
$ws.Range("H43").Value = 65672.336
$ws.Range("I43").Value = 6996
$ws.Range("J43").Value = 95010.5
$ws.Range("K43").Value = 6996
$ws.Range("L43").Value = 95010.5
$ws.Range("M43").Value = -6803
$ws.Range("N43").Value = -95396.5

# LTW!row46 - Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 2870
$ws.Range("J46").Value = 3111.111
$ws.Range("L46").Value = 3111.111
$ws.Range("N46").Value = -3487.111

# LTW!row126 - Battered Books / Saiga Leather
$ws.Range("H126").Value = 2109.0908
$ws.Range("I126").Value = 2316
$ws.Range("J126").Value = 1936.6666
$ws.Range("K126").Value = 6948
$ws.Range("L126").Value = 5809.9998
$ws.Range("M126").Value = -4478
$ws.Range("N126").Value = -10749.9998

# LTW!row141 - Just Generally Freezing / Gargantuaskin Trousers of Striking
$ws.Range("H141").Value = 29900
$ws.Range("J141").Value = 29900
$ws.Range("L141").Value = 29900
$ws.Range("N141").Value = -40260

$ws = $wb.Worksheets.Item("WVR")
# WVR!row19 - Dirt Cheap / Stablehand's Hat
$ws.Range("H19").Value = 9361
$ws.Range("I19").Value = 5000
$ws.Range("J19").Value = 10451.25
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 10451.25
$ws.Range("M19").Value = -4826
$ws.Range("N19").Value = -10799.25

# WVR!row29 - Getting Handsy / Cotton Dress Gloves
$ws.Range("H29").Value = 1811
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 1811
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 1811
$ws.Range("N29").Value = -2391
$ws.Range("M29").ClearContents()

# WVR!row70 - An Account of My Boots / Holy Rainbow Shoes
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# WVR!row73 - Soot in My Hair and Scars on My Feet (L) / Holy Rainbow Shoes
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
